$d = $word.ActiveDocument

# --- Locate the anchor paragraph: the bullet ending in
# "...counts from 1 - 10, 100, 1000" (last child bullet under
# "Breaking the Problem Apart"). We resolve it to a paragraph *index*
# (rather than holding on to the Range object across edits) because
# paragraph ranges/list levels need to be re-queried fresh after each
# insertion.
$rng = $d.Content
$found = $rng.Find.Execute("The goals is to know what finger she lands on when she counts from 1")
if (-not $found) {
    throw "Could not find anchor paragraph text"
}
[void]$rng.Expand(4)   # wdParagraph - expand the found hit to its whole paragraph
$anchorStart = $rng.Start

$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $anchorStart) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq 0) {
    throw "Could not resolve anchor paragraph index"
}

# --- Insert new bullet #1: "Identify potential Solutions" at ilvl 0 ---
$anchorPara = $d.Paragraphs.Item($anchorIndex)
$ip = $anchorPara.Range
$ip.Collapse(0)   # wdCollapseEnd
$ip.InsertParagraphAfter()

$p1 = $d.Paragraphs.Item($anchorIndex + 1)
$p1.Range.Text = "Identify potential Solutions"
$p1.Range.ListFormat.ListOutdent()

# --- Insert new bullet #2: the factor-of-nine explanation at ilvl 1 ---
$p1Again = $d.Paragraphs.Item($anchorIndex + 1)
$ip2 = $p1Again.Range
$ip2.Collapse(0)
$ip2.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item($anchorIndex + 2)
$p2.Range.Text = "If we associate the thumb as always being a factor of nine, we can then divide the number by nine (10/9) = 1.1111, then take the 10-(9*1) = 1 and that is how many fingers you would move from the thumb."
$p2.Range.ListFormat.ListIndent()
